# Auto-generated Excel COM-interop script
# Applies market-price / profit-column updates across all Leve sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1116.6666
$ws.Range("I18").Value = 1116.6666
$ws.Range("K18").Value = 1116.6666
$ws.Range("M18").Value = -832.6666
$ws.Range("H31").Value = 150
$ws.Range("I31").Value = 150
$ws.Range("K31").Value = 450
$ws.Range("M31").Value = -220
$ws.Range("H34").Value = 4360
$ws.Range("I34").Value = 3933.6667
$ws.Range("K34").Value = 3933.6667
$ws.Range("M34").Value = -3730.6667
$ws.Range("H36").Value = 4360
$ws.Range("I36").Value = 3933.6667
$ws.Range("K36").Value = 3933.6667
$ws.Range("M36").Value = -3218.6667
$ws.Range("H64").Value = 1450.25
$ws.Range("J64").Value = 1450.25
$ws.Range("L64").Value = 1450.25
$ws.Range("N64").Value = -1946.25
$ws.Range("H67").Value = 1450.25
$ws.Range("J67").Value = 1450.25
$ws.Range("L67").Value = 1450.25
$ws.Range("N67").Value = -3166.25
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 3333.3333
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -10384
$ws.Range("H138").Value = 4140.5454
$ws.Range("J138").Value = 5978.5713
$ws.Range("L138").Value = 17935.7139
$ws.Range("N138").Value = -28215.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 699
$ws.Range("I2").Value = 699
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 699
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -586
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 30
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 30
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 82
$ws.Range("N5").ClearContents()
$ws.Range("H45").Value = 1700
$ws.Range("I45").Value = 1400
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1023
$ws.Range("N45").Value = -2754
$ws.Range("H101").Value = 173351
$ws.Range("J101").Value = 173351
$ws.Range("L101").Value = 173351
$ws.Range("N101").Value = -179841
$ws.Range("H116").Value = 699
$ws.Range("I116").Value = 699
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 699
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1595
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 699
$ws.Range("I3").Value = 699
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 699
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -585
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 85
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 565
$ws.Range("J22").Value = 401
$ws.Range("L22").Value = 401
$ws.Range("N22").Value = -747
$ws.Range("H33").Value = 4709
$ws.Range("J33").Value = 4841.3335
$ws.Range("L33").Value = 4841.3335
$ws.Range("N33").Value = -5513.3335
$ws.Range("H81").Value = 71998.75
$ws.Range("J81").Value = 71998.75
$ws.Range("L81").Value = 71998.75
$ws.Range("N81").Value = -74120.75
$ws.Range("H84").Value = 71998.75
$ws.Range("J84").Value = 71998.75
$ws.Range("L84").Value = 215996.25
$ws.Range("N84").Value = -226604.25
$ws.Range("H135").Value = 49995
$ws.Range("J135").Value = 49995
$ws.Range("L135").Value = 49995
$ws.Range("N135").Value = -60135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11429543
$ws.Range("I6").Value = 11429543
$ws.Range("K6").Value = 11429543
$ws.Range("M6").Value = -11429430
$ws.Range("H16").Value = 1082.6
$ws.Range("I16").Value = 1133.3334
$ws.Range("J16").Value = 1006.5
$ws.Range("K16").Value = 1133.3334
$ws.Range("L16").Value = 1006.5
$ws.Range("M16").Value = -846.3334
$ws.Range("N16").Value = -1580.5
$ws.Range("H17").Value = 3101.6
$ws.Range("I17").Value = 1004
$ws.Range("J17").Value = 4500
$ws.Range("K17").Value = 1004
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = -830
$ws.Range("N17").Value = -4848
$ws.Range("H22").Value = 583
$ws.Range("I22").Value = 585
$ws.Range("J22").Value = 575
$ws.Range("K22").Value = 585
$ws.Range("L22").Value = 575
$ws.Range("M22").Value = -235
$ws.Range("N22").Value = -1275
$ws.Range("H28").Value = 45199.4
$ws.Range("J28").Value = 45199.4
$ws.Range("L28").Value = 45199.4
$ws.Range("N28").Value = -45689.4
$ws.Range("H31").Value = 9854.5625
$ws.Range("I31").Value = 5585.5713
$ws.Range("J31").Value = 13174.889
$ws.Range("K31").Value = 5585.5713
$ws.Range("L31").Value = 13174.889
$ws.Range("M31").Value = -5290.5713
$ws.Range("N31").Value = -13764.889
$ws.Range("H34").Value = 9854.5625
$ws.Range("I34").Value = 5585.5713
$ws.Range("J34").Value = 13174.889
$ws.Range("K34").Value = 5585.5713
$ws.Range("L34").Value = 13174.889
$ws.Range("M34").Value = -5383.5713
$ws.Range("N34").Value = -13578.889
$ws.Range("H36").Value = 4999.5
$ws.Range("J36").Value = 3498
$ws.Range("L36").Value = 3498
$ws.Range("N36").Value = -4274
$ws.Range("H40").Value = 4999.5
$ws.Range("J40").Value = 3498
$ws.Range("L40").Value = 3498
$ws.Range("N40").Value = -3818
$ws.Range("H44").Value = 29999
$ws.Range("I44").Value = 29998.334
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 29998.334
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -29556.334
$ws.Range("N44").Value = -30884
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H94").Value = 1123.5
$ws.Range("J94").Value = 1799.6
$ws.Range("L94").Value = 1799.6
$ws.Range("N94").Value = -2701.6
$ws.Range("H99").Value = 1432714.2
$ws.Range("I99").Value = 1668000
$ws.Range("K99").Value = 1668000
$ws.Range("M99").Value = -1666502
$ws.Range("H113").Value = 1082.6
$ws.Range("I113").Value = 1133.3334
$ws.Range("J113").Value = 1006.5
$ws.Range("K113").Value = 1133.3334
$ws.Range("L113").Value = 1006.5
$ws.Range("M113").Value = 1036.6666
$ws.Range("N113").Value = -5346.5
$ws.Range("H126").Value = 1432714.2
$ws.Range("I126").Value = 1668000
$ws.Range("K126").Value = 5004000
$ws.Range("M126").Value = -5001530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 3373.75
$ws.Range("I14").Value = 3373.75
$ws.Range("K14").Value = 10121.25
$ws.Range("M14").Value = -9948.25
$ws.Range("H38").Value = 902.1111
$ws.Range("I38").Value = 40
$ws.Range("J38").Value = 1148.4286
$ws.Range("K38").Value = 120
$ws.Range("L38").Value = 3445.2858
$ws.Range("M38").Value = 227
$ws.Range("N38").Value = -4139.2858
$ws.Range("H52").Value = 1750
$ws.Range("J52").Value = 1750
$ws.Range("L52").Value = 5250
$ws.Range("N52").Value = -5782
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1966.8334
$ws.Range("I16").Value = 1966.8334
$ws.Range("K16").Value = 1966.8334
$ws.Range("M16").Value = -1796.8334
$ws.Range("H22").Value = 2713.8572
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1205
$ws.Range("H27").Value = 2713.8572
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1393
$ws.Range("H46").Value = 403518.8
$ws.Range("I46").Value = 1001000
$ws.Range("K46").Value = 1001000
$ws.Range("M46").Value = -1000812

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 26960.375
$ws.Range("I51").Value = 28510.334
$ws.Range("J51").Value = 26030.4
$ws.Range("K51").Value = 28510.334
$ws.Range("L51").Value = 26030.4
$ws.Range("M51").Value = -28000.334
$ws.Range("N51").Value = -27050.4
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

Write-Host "Applied 233 cell updates and 13 cell clears across Leve profit sheets."
